$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Individuals": insert a new "alive" column before the "age" column
# (shifting age/spouse/death one column to the right) and populate it with a
# boolean telling whether the individual has no recorded death date.
# ---------------------------------------------------------------------------
$individuals = $wb.Worksheets.Item("Individuals")

# Shifts existing G:I ("age","spouse","death") -> H:J, carrying formatting
# (the new column G inherits the bold/border header style from its neighbour).
$individuals.Columns("G:G").Insert()
$individuals.Range("G1").Value = "alive"

$aliveByRow = @{
    2 = $true;  3 = $true;  4 = $true;  5 = $true
    6 = $true;  7 = $true;  8 = $true;  9 = $true
    10 = $false; 11 = $false; 12 = $false
}
foreach ($row in $aliveByRow.Keys) {
    $individuals.Cells.Item($row, 7).Value = $aliveByRow[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Families": append a new "are divorced" column after "divorced".
# ---------------------------------------------------------------------------
$families = $wb.Worksheets.Item("Families")

# Copy the bold/border header formatting from the neighbouring header cell
# so the new header reuses the existing style instead of creating a new one.
$families.Range("I1").Copy()
$families.Range("J1").PasteSpecial(-4122)
$families.Range("J1").Value = "are divorced"

$families.Range("J5").Value = $false
